# Update workbook for new battery: Turnigy 5000 mAh 6S 25C LiPo
$wb = $excel.ActiveWorkbook

# --- "weight analysis" sheet ---
$ws = $wb.Worksheets.Item("weight analysis")

# C17 (battery capacity-based energy term) becomes a formula based on new battery specs
$ws.Range("C17").Formula = "=5000*3.6*22.2/C26"

# C25 (thrust/weight? total weight force) - new combined mass constant
$ws.Range("C25").Formula = "=14.031*9.8"

# C26 (battery weight force) becomes a formula: 0.655 kg * 9.81
$ws.Range("C26").Formula = "=0.655*9.81"

# C31 changes from a formula to a hardcoded constant
$ws.Range("C31").Value = 0.511

# --- "PropWGlider" sheet ---
$ws2 = $wb.Worksheets.Item("PropWGlider")
$ws2.Range("J19").Value = 19.4
$ws2.Range("B20").Value = 0.41
$ws2.Range("B22").Value = 0.066
$ws2.Range("B24").Value = 0.037

# --- "PropWOGlider" sheet ---
$ws3 = $wb.Worksheets.Item("PropWOGlider")
$ws3.Range("J19").Value = 18.75
$ws3.Range("B20").Value = 0.425
$ws3.Range("B22").Value = 0.061

# Recalculate all formulas so cached/dependent values are refreshed
$excel.CalculateFullRebuild()

# --- Selection / view state updates (cosmetic, matches diff) ---
$ws.Application.Goto($ws.Range("C31"), $true)
$ws.Range("C31").Select() | Out-Null

$ws2.Range("J27").Select() | Out-Null
